$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A50: "03-11-2025" -------------------------------------------------
# Plain Value assignment would make Excel's input-parser treat a
# date-shaped string as an actual date (numeric serial + date number
# format). The source file stores every date in this column as literal
# text, so we build the text via a formula (whose *result* type is
# string, not user input) on a scratch cell, then copy/paste-special
# just the value onto the target cell - this keeps it a literal string.
$ws.Cells.Item(50, 1).Borders.LineStyle = 1
$ws.Cells.Item(1000, 1).Formula = "=""03-11-2025"""
$ws.Cells.Item(1000, 1).Copy()
$ws.Cells.Item(50, 1).PasteSpecial(-4163)
$ws.Cells.Item(1000, 1).EntireRow.Delete()

# --- B50: gold-price text ----------------------------------------------
$ws.Cells.Item(50, 2).Value = "The price of gold in India today is ₹12,317 per gram for 24 karat gold, ₹11,290 per gram for 22 karat gold and ₹9,238 per gram for 18 karat gold (also called 999 gold)."
